$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the student record in row 2 with the new cadastro data.
$ws.Range("B2").Value = "Teste"

# RA and "Registro de Matricula" are identifier-like strings; keep them as
# text even though they look numeric (force text storage, then drop back
# to the default/general style so no stray formatting is left behind).
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "1213"
$ws.Range("C2").NumberFormat = "General"
$ws.Range("C2").Style = "Normal"

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "23"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").Value = "2025-02-24 00:41:38"

# New "Registro de Matricula" column needs its own width now that it holds data.
$ws.Columns.Item(4).ColumnWidth = 17.75

# Reflect the current selection used while reviewing the new row.
$ws.Range("A2:E4").Select() | Out-Null
